# Insert a new price-record row for "Angeleno / Segunda" (2022-03-17) above
# the existing row 87, pushing the following records (old rows 87-103) down
# by one row to rows 88-104. This mirrors the weekly Fruta/Hortalizas data
# refresh described in the commit message "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing row 87 (and everything below it) down one row.
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new record.
$ws.Cells.Item(87, 1).Value  = 5
$ws.Cells.Item(87, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(87, 3).Value  = "Maule"
$ws.Cells.Item(87, 4).Value  = 44637
$ws.Cells.Item(87, 5).Value  = 7
$ws.Cells.Item(87, 6).Value  = "Fruta"
$ws.Cells.Item(87, 7).Value  = 100103
$ws.Cells.Item(87, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(87, 9).Value  = 100103002
$ws.Cells.Item(87, 10).Value = "Ciruela"
$ws.Cells.Item(87, 11).Value = "Angeleno"
$ws.Cells.Item(87, 12).Value = "Segunda"
$ws.Cells.Item(87, 13).Value = 200
$ws.Cells.Item(87, 14).Value = 6000
$ws.Cells.Item(87, 15).Value = 6000
$ws.Cells.Item(87, 16).Value = 6000
$ws.Cells.Item(87, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(87, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(87, 19).Value = 333
$ws.Cells.Item(87, 20).Value = 18
